$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.758.73"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.618.41"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.15"
$ws.Range("E5").Value = "  +1.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.49"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.632.68"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  +3.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("E13").Value = "  +1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.075.49"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.775.64"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.74"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.626.88"

$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.68"
$ws.Range("E20").Value = "  +3.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.67"
$ws.Range("E21").Value = "  +2.41%  "

$ws.Range("E22").Value = "  +0.84%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.91"
$ws.Range("E24").Value = "  +1.29%  "

$ws.Range("E25").Value = "  +0.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.734.17"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0847"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.50"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.59"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.96"
$ws.Range("E34").Value = "  +4.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "151.65"
$ws.Range("E35").Value = "  -3.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.03"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.888"
$ws.Range("E38").Value = "  +6.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.855"
$ws.Range("E39").Value = "  +0.99%  "

$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "293.33"
$ws.Range("E43").Value = "  -4.28%  "

$ws.Range("E44").Value = "  +0.40%  "

$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("E46").Value = "  -2.38%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.82"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.98"
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.31"
$ws.Range("E51").Value = "  +0.19%  "
